# Add a new September transaction entry at the top of the "2024" sheet's
# September_Details / September_Date columns (R/S), pushing all existing
# rows (47-192) down by one row. This also shifts the "Broadband" group
# label from A192 down to A193.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 47 (the first September data row),
# shifting everything below it (including the trailing "Broadband" row)
# down by one.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row with the new transaction.
$ws.Range("R47").Value = "money google"
$ws.Range("S47").Value = "2024-09-22 20:17:45"
